$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs updating from
# 2023-09-03 (45172) to 2023-09-06 (45175) for every data row (rows 2-269).
$startRow = 2
$endRow = 269
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(45175)

for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
